# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Pais")

# Update the "last updated" timestamp string (row 1, column A)
$ws.Range("A1").Value = "Datos actualizados a 13 de Junio de 2020 a las 23:42"

# Estados Unidos (row 4)
$ws.Range("B4").Value = 2138947
$ws.Range("C4").Value = 22025
$ws.Range("D4").Value = 849107
$ws.Range("E4").Value = 1172414
$ws.Range("G4").Value = 601
$ws.Range("H4").Value = 117426

# Brasil (row 5)
$ws.Range("B5").Value = 850514
$ws.Range("C5").Value = 20612
$ws.Range("E5").Value = 380184
$ws.Range("G5").Value = 819
$ws.Range("H5").Value = 42720

# Alemania (row 12)
$ws.Range("B12").Value = 187423
$ws.Range("C12").Value = 172
$ws.Range("E12").Value = 6656

# Bulgaria (row 89)
$ws.Range("B89").Value = 3266
$ws.Range("C89").Value = 75
$ws.Range("D89").Value = 1723
$ws.Range("E89").Value = 1371

# Row 145 was "Malaui" -> becomes "Togo" with updated figures
$ws.Range("A145").Value = "Togo"
$ws.Range("B145").Value = 530
$ws.Range("C145").Value = 5
$ws.Range("D145").Value = 291
$ws.Range("E145").Value = 226
$ws.Range("G145").Value = 0
$ws.Range("H145").Value = 13

# Row 146 was "Togo" -> becomes "Malaui" with the previous Malaui figures
$ws.Range("A146").Value = "Malaui"
$ws.Range("B146").Value = 529
$ws.Range("C146").Value = 48
$ws.Range("D146").Value = 66
$ws.Range("E146").Value = 458
$ws.Range("G146").Value = 1
$ws.Range("H146").Value = 5

# Gambia (row 192)
$ws.Range("D192").Value = 24
$ws.Range("E192").Value = 3

# Row 210 was "Montserrat" -> becomes "Seychelles"
$ws.Range("A210").Value = "Seychelles"
$ws.Range("D210").Value = 11
$ws.Range("H210").Value = 0

# Row 211 was "Seychelles" -> becomes "Montserrat"
$ws.Range("A211").Value = "Montserrat"
$ws.Range("D211").Value = 10
$ws.Range("H211").Value = 1
